$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 185.42857
$ws.Range("I33").Value = 141.33333
$ws.Range("K33").Value = 141.33333
$ws.Range("M33").Value = 87.66667000000001
# Row 70
$ws.Range("H70").Value = 3499.5715
$ws.Range("I70").Value = 3499.5715
$ws.Range("K70").Value = 10498.7145
$ws.Range("M70").Value = -10228.7145
# Row 73
$ws.Range("H73").Value = 3499.5715
$ws.Range("I73").Value = 3499.5715
$ws.Range("K73").Value = 10498.7145
$ws.Range("M73").Value = -9562.7145
# Row 98
$ws.Range("H98").Value = 9099.583000000001
$ws.Range("I98").Value = 6720
$ws.Range("J98").Value = 20997.5
$ws.Range("K98").Value = 6720
$ws.Range("L98").Value = 20997.5
$ws.Range("M98").Value = -5222
$ws.Range("N98").Value = -23993.5
# Row 116
$ws.Range("H116").Value = 3997.5
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
# Row 122
$ws.Range("H122").Value = 9099.583000000001
$ws.Range("I122").Value = 6720
$ws.Range("J122").Value = 20997.5
$ws.Range("K122").Value = 20160
$ws.Range("L122").Value = 62992.5
$ws.Range("M122").Value = -17710
$ws.Range("N122").Value = -67892.5

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 549.375
$ws.Range("I32").Value = 549.375
$ws.Range("K32").Value = 549.375
$ws.Range("M32").Value = -262.375
# Row 104
$ws.Range("H104").Value = 39000
$ws.Range("J104").Value = 39000
$ws.Range("L104").Value = 39000
$ws.Range("N104").Value = -45988
# Row 118
$ws.Range("H118").Value = 1100000
$ws.Range("I118").Value = 1100000
$ws.Range("K118").Value = 1100000
$ws.Range("M118").Value = -1098343
# Row 130
$ws.Range("H130").Value = 74000
$ws.Range("J130").Value = 74000
$ws.Range("L130").Value = 74000
$ws.Range("N130").Value = -84040
# Row 131
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 16
$ws.Range("H16").Value = 1942
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 1942
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 1942
$ws.Range("N16").Value = -2282
$ws.Range("M16").ClearContents()
# Row 22
$ws.Range("H22").Value = 916.5
$ws.Range("I22").Value = 859.8
$ws.Range("K22").Value = 859.8
$ws.Range("M22").Value = -686.8
# Row 88
$ws.Range("H88").Value = 20287.2
$ws.Range("J88").Value = 20287.2
$ws.Range("L88").Value = 20287.2
$ws.Range("N88").Value = -21099.2
# Row 91
$ws.Range("H91").Value = 20287.2
$ws.Range("J91").Value = 20287.2
$ws.Range("L91").Value = 20287.2
$ws.Range("N91").Value = -23095.2
# Row 106
$ws.Range("H106").Value = 4075
$ws.Range("J106").Value = 4075
$ws.Range("L106").Value = 4075
$ws.Range("N106").Value = -6599

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 33
$ws.Range("H33").Value = 38099.234
$ws.Range("I33").Value = 8526.857
$ws.Range("K33").Value = 8526.857
$ws.Range("M33").Value = -8147.857
# Row 132
$ws.Range("H132").Value = 2986.8
$ws.Range("I132").Value = 2988.5
$ws.Range("J132").Value = 2980
$ws.Range("K132").Value = 8965.5
$ws.Range("L132").Value = 8940
$ws.Range("M132").Value = -6435.5
$ws.Range("N132").Value = -14000
# Row 140
$ws.Range("H140").Value = 62999
$ws.Range("J140").Value = 62999
$ws.Range("L140").Value = 62999
$ws.Range("N140").Value = -73359

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 80
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
# Row 81
$ws.Range("H81").Value = 49999.168
$ws.Range("J81").Value = 49999.168
$ws.Range("L81").Value = 149997.504
$ws.Range("N81").Value = -152243.504
# Row 83
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
# Row 84
$ws.Range("H84").Value = 49999.168
$ws.Range("J84").Value = 49999.168
$ws.Range("L84").Value = 449992.512
$ws.Range("N84").Value = -461224.512
# Row 115
$ws.Range("H115").Value = 1000
$ws.Range("I115").Value = 1000
$ws.Range("K115").Value = 3000
$ws.Range("M115").Value = -1825
# Row 131
$ws.Range("H131").Value = 615
$ws.Range("I131").Value = 620
$ws.Range("J131").Value = 600
$ws.Range("K131").Value = 1860
$ws.Range("L131").Value = 1800
$ws.Range("M131").Value = 3180
$ws.Range("N131").Value = -11880

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Range("H11").Value = 26125000
$ws.Range("I11").Value = 26125000
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 26125000
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -26124861
$ws.Range("N11").ClearContents()
# Row 12
$ws.Range("H12").Value = 1501.5
$ws.Range("I12").Value = 2003
$ws.Range("J12").Value = 1000
$ws.Range("K12").Value = 2003
$ws.Range("L12").Value = 1000
$ws.Range("M12").Value = -1863
$ws.Range("N12").Value = -1280

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 933
$ws.Range("I22").Value = 800
$ws.Range("J22").Value = 999.5
$ws.Range("K22").Value = 800
$ws.Range("L22").Value = 999.5
$ws.Range("M22").Value = -505
$ws.Range("N22").Value = -1589.5
# Row 27
$ws.Range("H27").Value = 933
$ws.Range("I27").Value = 800
$ws.Range("J27").Value = 999.5
$ws.Range("K27").Value = 800
$ws.Range("L27").Value = 999.5
$ws.Range("M27").Value = -693
$ws.Range("N27").Value = -1213.5
# Row 38
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
# Row 122
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 3173.8333
$ws.Range("I81").Value = 3173.8333
$ws.Range("K81").Value = 6347.6666
$ws.Range("M81").Value = -5286.6666
# Row 84
$ws.Range("H84").Value = 3173.8333
$ws.Range("I84").Value = 3173.8333
$ws.Range("K84").Value = 31738.333
$ws.Range("M84").Value = -26434.333
# Row 122
$ws.Range("H122").Value = 1611.48
$ws.Range("I122").Value = 1621.174
$ws.Range("K122").Value = 4863.522
$ws.Range("M122").Value = -2413.522
